$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Ticketart attribute definition: square-bracket Dropdown[...] -> Dropdown(...)
$ws.Range("D3").Value = "Ticketart:Dropdown(Bus,Zug,U-Bahn); Häufigkeit:Dropdown(Täglich,Wöchentlich,Selten)"

# Scroll the sheet view so row 2 is the top-left visible row
$ws.Application.ActiveWindow.ScrollRow = 2

# Reset row 8's height back to the sheet's default (remove explicit custom height)
$ws.Rows.Item(8).AutoFit()
